$wb = $excel.ActiveWorkbook

# Sheet "展览": F3 1250 -> 1251, F4 2741 -> 2747
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1251
$ws1.Range("F4").Value = 2747

# Sheet "全部类型": F5 1250 -> 1251, F6 2741 -> 2747
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1251
$ws4.Range("F6").Value = 2747
